# feat: fill required columns
#
# Column G = "收件者姓名 / Tên người nhận" (receiver name)
# Column I = "收件電話 / Số điện thoại nhận hàng" (receiver phone)
# Several data rows are missing these required values; fill the blanks
# with the placeholder "None" so every row has a value in the required
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "receiver phone" (column I) cell is currently empty.
$rowsMissingPhone = @(4,5,6,7,8,9,10,11,12,13,14,15,16,17,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,37,38,39,40,41,42,43,44,46,47)

foreach ($r in $rowsMissingPhone) {
    $ws.Cells.Item($r, 9).Value = "None"
}

# Row 39 is also missing the "receiver name" (column G) value.
$ws.Cells.Item(39, 7).Value = "None"
